# "Updated symbol list" data refresh (GitHub Actions bot, 2023-01-26 10:58 UTC):
# refreshed Price (col D) / Volume(1h) (col E) figures for every coin row, and for
# rows 7-17 the coin ranking reshuffled so the Coin (B) / Link (C) in each of those
# rows now refers to a different coin as well.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cellNames = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "B7", "C7", "D7", "E7", "B8", "C8",
    "D8", "E8", "B9", "C9", "D9", "E9", "B10", "C10", "D10", "E10", "B11", "C11", "D11", "E11",
    "B12", "C12", "D12", "E12", "B13", "C13", "D13", "E13", "B14", "C14", "D14", "E14", "B15",
    "C15", "D15", "E15", "B16", "C16", "D16", "E16", "B17", "C17", "D17", "E17", "E18", "D19",
    "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "E24", "D25", "E25", "D26",
    "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43",
    "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50",
    "E50", "D51", "E51"
)

# Columns D (Price) and E (Volume(1h)) hold numeric-looking text, e.g. "304.39" or
# "1.05%". A leading apostrophe is Excel's text-qualifier: it forces the assignment
# to land as a literal text value (matching the original cells) instead of being
# auto-converted into a real number/percentage.
$cellValues = @(
    "'304.39", "'1.05%", "'35.86", "'2.05%", "'5.086", "'0.95%", "'0.08038", "'1.14%", "'1.915",
    "'0.25%", "KuCoinToken", "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs",
    "'7.732", "'-0.92%", "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx",
    "'0.9287", "'0.72%", "LiechtensteinCryptoassetsExchange",
    "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx", "'0.1374",
    "'4.11%", "WazirX", "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx", "'0.1891", "'2.68%",
    "MandalaExchangeToken", "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx",
    "'0.09093", "'-4.78%", "BitrueCoin", "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr",
    "'0.03601", "'0.83%", "BitMartToken",
    "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx", "'0.09818", "'-0.21%",
    "BitForexToken", "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf", "'0.001447",
    "'3.25%", "TigerCash", "https://coinranking.com/coin/6hIn06L2+tigercash-tch", "'0.005921",
    "'2.10%", "LEO", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo", "'3.555", "'1.45%",
    "GateToken", "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt", "'4.183", "'3.73%",
    "'1.81%", "'0.3458", "'0.91%", "'0.1333", "'1.75%", "'4.903", "'-3.07%", "'0.2514", "'2.14%",
    "'0.04462", "'-0.88%", "'0.56%", "'0.004782", "'-0.24%", "'0.0001484", "'18.66%", "'0.0003140",
    "'4.63%", "'0.01951", "'4.16%", "'0.04874", "'3.45%", "'0.007625", "'1.90%", "'0.009188",
    "'-5.79%", "'0.1371", "'3.65%", "'0.002105", "'-0.25%", "'0.01137", "'18.59%", "'0.00006408",
    "'3.70%", "'0.00000000752", "'0.20%", "'64.67", "'0.29%", "'0.001195", "'-19.79%",
    "'0.00002105", "'0.20%", "'0.0002005", "'0.20%"
)

for ($i = 0; $i -lt $cellNames.Length; $i++) {
    $rng = $ws.Range($cellNames[$i])
    $rng.Value = $cellValues[$i]
    # Assigning a leading-apostrophe string marks the cell with Excel's "quote
    # prefix" (number-stored-as-text) style flag. Clear the format right back so the
    # cell keeps the plain/default style it had before (these cells carry no other
    # formatting), matching the original text cells exactly.
    $rng.ClearFormats()
}
